$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 170822
$ws.Range("C4").Value = 161645
$ws.Range("C7").Value = 5.37
$ws.Range("C8").Value = 65.79000000000001
